# Update the "Förändrad" (Changed) date column (C) for rows 2-18
# from 2023-10-04 (45203) to 2023-10-05 (45204), keeping formatting intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 18; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45203) {
        $cell.Value = 45204
    }
}
